$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.816.53"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.57%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.502.97"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -3.59%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.54"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -4.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "192.68"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -4.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.612"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -2.42%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.493.14"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.53%  "
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.205"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -6.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "51.56"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -4.46%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000287"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -6.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.16"
$ws.Range("D14").ClearFormats()
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.059.89"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -3.48%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "644.96"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -4.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.740.11"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.499.35"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -4.89%  "
$ws.Range("E19").Value = "  -3.66%  "
$ws.Range("E20").Value = "  -1.77%  "
$ws.Range("E21").Value = "  -3.38%  "
$ws.Range("E22").Value = "  -4.92%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "18.04"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.76%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.34"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.95%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "99.21"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -5.49%  "
$ws.Range("E26").Value = "  -7.34%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.89"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -4.36%  "
$ws.Range("E28").Value = "  -4.24%  "
$ws.Range("E29").Value = "  -4.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.74"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -4.30%  "
$ws.Range("E31").Value = "  -9.48%  "
$ws.Range("E32").Value = "  -6.41%  "
$ws.Range("E33").Value = "  -4.35%  "
$ws.Range("E34").Value = "  -4.72%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "61.54"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.78%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "564.69"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +9.58%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.13"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +61.76%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.717.24"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -5.86%  "
$ws.Range("E39").Value = "  +0.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0791"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -8.93%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.60"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.58%  "
$ws.Range("E42").Value = "  -3.79%  "
$ws.Range("E43").Value = "  -3.72%  "
$ws.Range("E44").Value = "  -1.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "34.37"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -5.89%  "
$ws.Range("E46").Value = "  -3.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.38"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.25%  "
$ws.Range("E48").Value = "  -6.39%  "
$ws.Range("E49").Value = "  -4.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.998"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.42%  "
$ws.Range("E51").Value = "  -5.08%  "
